{"js": "// Replace the placeholder citation codes (Ref-XXXXXX) throughout the\n// document. Every citation marker inside a given paragraph collapses to a\n// single new code for that paragraph (several distinct old codes -> one new\n// code), per the commit's \"citation check\" dataset update.\nconst paragraphRefMap = {\n  1: \"Ref-f011992\",\n  3: \"Ref-u116920\",\n  5: \"Ref-s600855\",\n  6: \"Ref-s418542\",\n  8: \"Ref-s305680\",\n  9: \"Ref-s366502\",\n  10: \"Ref-u823257\",\n};\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst refPattern = /Ref-[A-Za-z0-9]+/g;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const newRef = paragraphRefMap[i];\n  if (!newRef) continue;\n\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text || \"\";\n  const oldRefs = new Set();\n  let m;\n  refPattern.lastIndex = 0;\n  while ((m = refPattern.exec(text)) !== null) {\n    oldRefs.add(m[0]);\n  }\n\n  for (const oldRef of oldRefs) {\n    if (oldRef === newRef) continue;\n    const results = paragraph.search(oldRef, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n    for (const range of results.items) {\n      range.insertText(newRef, \"Replace\");\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Replace the placeholder citation codes (Ref-XXXXXX) throughout the\n# document. Every citation marker inside a given paragraph collapses to a\n# single new code for that paragraph (several distinct old codes -> one new\n# code), per the commit's \"citation check\" dataset update.\n\n$d = $word.ActiveDocument\n\n$paragraphRefMap = @{\n    2  = \"Ref-f011992\"\n    4  = \"Ref-u116920\"\n    6  = \"Ref-s600855\"\n    7  = \"Ref-s418542\"\n    9  = \"Ref-s305680\"\n    10 = \"Ref-s366502\"\n    11 = \"Ref-u823257\"\n}\n\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    if (-not $paragraphRefMap.ContainsKey($i)) { continue }\n    $newRef = $paragraphRefMap[$i]\n\n    $paragraph = $d.Paragraphs.Item($i)\n    $text = $paragraph.Range.Text\n\n    $oldRefs = New-Object System.Collections.Generic.HashSet[string]\n    $regexMatches = [regex]::Matches($text, \"Ref-[A-Za-z0-9]+\")\n    foreach ($m in $regexMatches) {\n        [void]$oldRefs.Add($m.Value)\n    }\n\n    foreach ($oldRef in $oldRefs) {\n        if ($oldRef -eq $newRef) { continue }\n        $searchRange = $d.Paragraphs.Item($i).Range\n        $searchRange.Find.Execute($oldRef, $false, $true, $false, $false, $false, $true, 1, $false, $newRef, 2)\n    }\n}\n\n$d.Save()\n"}
